$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inventory for check stock")
$ws.Range("F2").Value = "01A000"
$ws.Range("F2").Select()
